$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the "import" cell: semicolons -> commas, drop the trailing semicolon.
$ws.Range("B2").Value = "com.blackknight.demo.models.MortgageRequest," + [char]10 + "com.blackknight.demo.models.Address," + [char]10 + "Com.blackknight.demo.models.Loan"

# 2. Collapse the rule table from 4 "MortgageRequest" condition columns (B:E)
#    down to 2 (B:C) by removing columns C and D entirely (this also shifts the
#    old F:I columns left by two, producing the new B:G, 7-column layout and
#    automatically fixing the dimension + the B6:E6/F6:G6 -> B6:C6/D6:E6 merges).
$ws.Range("C1:D9").EntireColumn.Delete()

# 3. After the shift, the old "annualIncome >= $param" sub-condition header
#    (now at C5) needs to go back to reading "CONDITON" (matching the other
#    sub-columns of the first condition group).
$ws.Range("C5").Value = "CONDITON"

# 4. The three separate sub-condition rules that used to live in B7, C7, D7
#    are combined into a single comma-joined string in B7.
$ws.Range("B7").Value = "`$address: mailingAddress,`$loan: loan,annualIncome >= `$param"

# 5. Restore the "Annual Income" column header (now in B8) with bold styling
#    to match the rest of the header row.
$ws.Range("B8").Value = "Annual Income"
$ws.Range("B8").Font.Bold = $true

# 6. The sample rule's first two data cells are now numeric values instead of
#    placeholder "x" text.
$ws.Range("B9").Value = 5000

# 7. Widen column B to fit the merged condition text.
$ws.Columns("B").ColumnWidth = 51.7

# 8. Restore the saved selection.
$ws.Range("D16").Select() | Out-Null
